$wb = $excel.ActiveWorkbook

# --- Sheet "مسکن ایلیا" (Ilia's housing finance tracker) ---
$ilia = $wb.Worksheets.Item("مسکن ایلیا")
$ilia.Activate()

# A new ledger entry was recorded on 20/4/1396: an ATM withdrawal of 200,000
# toman, partly (150) for a 30kg bag of rice. It occupies the first
# previously-blank row (131) right after the last populated row (130), and
# the "months" counter on the prior row is adjusted from 2 to 1.
$ilia.Range("C130").Value = 1

$ilia.Range("A131").Value = "20/4/1396"
$ilia.Range("B131").Value = -200000
$ilia.Range("C131").Value = 2
$ilia.Range("G131").Value = "از عابربانک گرفتم، 150 بابت بخشی از 30 کیلو برنج"

# Move the selection to the next blank row, mirroring where the user left
# off after data entry.
$ilia.Range("G132").Select()

# --- Sheet "برنامه 5 ساله" (5-year plan) ---
$plan = $wb.Worksheets.Item("برنامه 5 ساله")

# Ali's cash on hand was updated from 31,000 to 41,000.
$plan.Range("K16").Value = 41000

# This becomes the active sheet/tab, matching where the user ended up.
$plan.Activate()
